# Actualización automática 2025-07-03 12:35:11
# Updates the "CUMPLIMIENTO MENSUAL" sheet: refreshed PRESUPUESTO / VENTA /
# POR CUMPLIR / CUMPLIMIENTO figures for several product groups, plus the
# TOTAL row, and narrows/widens a few columns to fit the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Column widths (D, E, F) ---
$ws.Columns.Item(4).ColumnWidth = 12.15   # stored width 14 -> 13
$ws.Columns.Item(5).ColumnWidth = 16.15   # stored width 24 -> 17
$ws.Columns.Item(6).ColumnWidth = 24.15   # stored width 24 -> 25

# --- Row 2: 240X120 PORCELANATO ---
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 5820
$ws.Range("F2").Value = 0

# --- Row 3: 240X80 PORCELANATO ---
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 13728
$ws.Range("F3").Value = 0

# --- Row 4: FREGADEROS DE COCINA ---
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 646
$ws.Range("F4").Value = 0

# --- Row 6: GRIFERIAS ---
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 106.82
$ws.Range("F6").Value = 0

# --- Row 7: INODOROS ---
$ws.Range("C7").Value = 3200
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 3200
$ws.Range("F7").Value = 0

# --- Row 8: LAVABOS ---
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1000
$ws.Range("F8").Value = 0

# --- Row 12: PANELES DECORATIVOS ---
$ws.Range("D12").Value = 36
$ws.Range("E12").Value = 314
$ws.Range("F12").Value = 0.1028571428571429

# --- Row 13: PANELES PU ---
$ws.Range("C13").Value = 130
$ws.Range("E13").Value = 130

# --- Row 14: PANELES PVC ---
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = 2860.94
$ws.Range("E14").Value = -2620.94
$ws.Range("F14").Value = 11.92058333333333

# --- Row 15: PIEDRA SINTERIZADA ---
$ws.Range("C15").Value = 20690
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 20690
$ws.Range("F15").Value = 0

# --- Row 16: PORCELANATO ---
$ws.Range("C16").Value = 54721.23
$ws.Range("D16").Value = 2059.14
$ws.Range("E16").Value = 52662.09
$ws.Range("F16").Value = 0.03762963661452785

# --- Row 17: PUERTAS DE SEGURIDAD ---
$ws.Range("D17").Value = 326.24
$ws.Range("E17").Value = 815.76
$ws.Range("F17").Value = 0.2856742556917689

# --- Row 19: TOTAL ---
$ws.Range("C19").Value = 105212.87
$ws.Range("D19").Value = 5282.32
$ws.Range("E19").Value = 99930.54999999999
$ws.Range("F19").Value = 0.05020602517543719

Write-Output "Edit applied"
